# abiotic factor data entry
# Rows 52-71 were mislabeled as "F2" chamber samples; correct them to "F1".
# Only the SampleName (column E) text changes; the dependent formula
# columns (F:J) are formulas and will recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 52; $row -le 71; $row++) {
    $cell = $ws.Range("E$row")
    $cell.Value2 = $cell.Value2 -replace '^F2-', 'F1-'
}

# Restore the selection/active cell as recorded after the edit.
$ws.Range("F62").Select()
